$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenMap")
$ws.Cells.Clear()

# Row 1
$ws.Range("A1").Value = 'vstacks_t1~'
$ws.Range("B1").Value = 'vstacks_t5~'
$ws.Range("C1").Value = 'vstacks_w2~'
$ws.Range("H1").Value = 'C'

# Row 2
$ws.Range("H2").Value = 'ngfs'
$ws.Range("I2").Value = 'timeslice'

# Row 4
$ws.Range("A4").Value = '~ScenMap'
$ws.Range("G4").Value = '~ScenG'

# Row 5 - headers
$ws.Range("A5").Value = 'Oname'
$ws.Range("B5").Value = 'Name'
$ws.Range("C5").Value = 'Desc'
$ws.Range("D5").Value = 'GW'
$ws.Range("G5").Value = 'Scen'
$ws.Range("H5").Formula = '="sg_"&H2'
$ws.Range("I5").Formula = '="sg_"&I2'

# ---- Block 1: rows 6-12 (3 days) ----
# Row 6: Delayed transition
$ws.Range("A6").Formula = '=$A$1&TEXT(N6,"0000")'
$ws.Range("B6").Formula = '=G6'
$ws.Range("G6").Formula = '=H6&P6'
$ws.Range("H6").Value = 'Delayed transition'
$ws.Range("I6").Value = '3 days'
$ws.Range("N6").Value = 1
$ws.Range("P6").Value = '_3d'

# Row 7: Net Zero 2050
$ws.Range("A7").Formula = '=$A$1&TEXT(N7,"0000")'
$ws.Range("B7").Formula = '=G7'
$ws.Range("G7").Formula = '=H7&P7'
$ws.Range("H7").Value = 'Net Zero 2050'
$ws.Range("I7").Value = '3 days'
$ws.Range("N7").Value = 2
$ws.Range("P7").Value = '_3d'

# Row 8: NDCs
$ws.Range("A8").Formula = '=$A$1&TEXT(N8,"0000")'
$ws.Range("B8").Formula = '=G8'
$ws.Range("G8").Formula = '=H8&P8'
$ws.Range("H8").Value = 'NDCs'
$ws.Range("I8").Value = '3 days'
$ws.Range("N8").Value = 3
$ws.Range("P8").Value = '_3d'

# Row 9: Below 2deg
$ws.Range("A9").Formula = '=$A$1&TEXT(N9,"0000")'
$ws.Range("B9").Formula = '=G9'
$ws.Range("G9").Formula = '=H9&P9'
$ws.Range("H9").Value = 'Below 2deg'
$ws.Range("I9").Value = '3 days'
$ws.Range("N9").Value = 4
$ws.Range("P9").Value = '_3d'

# Row 10: Current Policies
$ws.Range("A10").Formula = '=$A$1&TEXT(N10,"0000")'
$ws.Range("B10").Formula = '=G10'
$ws.Range("G10").Formula = '=H10&P10'
$ws.Range("H10").Value = 'Current Policies'
$ws.Range("I10").Value = '3 days'
$ws.Range("N10").Value = 5
$ws.Range("P10").Value = '_3d'

# Row 11: Low demand
$ws.Range("A11").Formula = '=$A$1&TEXT(N11,"0000")'
$ws.Range("B11").Formula = '=G11'
$ws.Range("G11").Formula = '=H11&P11'
$ws.Range("H11").Value = 'Low demand'
$ws.Range("I11").Value = '3 days'
$ws.Range("N11").Value = 6
$ws.Range("P11").Value = '_3d'

# Row 12: Fragmented World
$ws.Range("A12").Formula = '=$A$1&TEXT(N12,"0000")'
$ws.Range("B12").Formula = '=G12'
$ws.Range("G12").Formula = '=H12&P12'
$ws.Range("H12").Value = 'Fragmented World'
$ws.Range("I12").Value = '3 days'
$ws.Range("N12").Value = 7
$ws.Range("P12").Value = '_3d'

# ---- Block 2: rows 13-19 (15 days) ----
# Row 13: Delayed transition
$ws.Range("A13").Formula = '=$B$1&TEXT(N13,"0000")'
$ws.Range("B13").Formula = '=G13'
$ws.Range("G13").Formula = '=H13&P13'
$ws.Range("H13").Formula = '=H6'
$ws.Range("I13").Value = '15 days'
$ws.Range("N13").Formula = '=N6'
$ws.Range("P13").Value = '_15d'

# Row 14: Net Zero 2050
$ws.Range("A14").Formula = '=$B$1&TEXT(N14,"0000")'
$ws.Range("B14").Formula = '=G14'
$ws.Range("G14").Formula = '=H14&P14'
$ws.Range("H14").Formula = '=H7'
$ws.Range("I14").Value = '15 days'
$ws.Range("N14").Formula = '=N7'
$ws.Range("P14").Value = '_15d'

# Row 15: NDCs
$ws.Range("A15").Formula = '=$B$1&TEXT(N15,"0000")'
$ws.Range("B15").Formula = '=G15'
$ws.Range("G15").Formula = '=H15&P15'
$ws.Range("H15").Formula = '=H8'
$ws.Range("I15").Value = '15 days'
$ws.Range("N15").Formula = '=N8'
$ws.Range("P15").Value = '_15d'

# Row 16: Below 2deg
$ws.Range("A16").Formula = '=$B$1&TEXT(N16,"0000")'
$ws.Range("B16").Formula = '=G16'
$ws.Range("G16").Formula = '=H16&P16'
$ws.Range("H16").Formula = '=H9'
$ws.Range("I16").Value = '15 days'
$ws.Range("N16").Formula = '=N9'
$ws.Range("P16").Value = '_15d'

# Row 17: Current Policies
$ws.Range("A17").Formula = '=$B$1&TEXT(N17,"0000")'
$ws.Range("B17").Formula = '=G17'
$ws.Range("G17").Formula = '=H17&P17'
$ws.Range("H17").Formula = '=H10'
$ws.Range("I17").Value = '15 days'
$ws.Range("N17").Formula = '=N10'
$ws.Range("P17").Value = '_15d'

# Row 18: Low demand
$ws.Range("A18").Formula = '=$B$1&TEXT(N18,"0000")'
$ws.Range("B18").Formula = '=G18'
$ws.Range("G18").Formula = '=H18&P18'
$ws.Range("H18").Formula = '=H11'
$ws.Range("I18").Value = '15 days'
$ws.Range("N18").Formula = '=N11'
$ws.Range("P18").Value = '_15d'

# Row 19: Fragmented World
$ws.Range("A19").Formula = '=$B$1&TEXT(N19,"0000")'
$ws.Range("B19").Formula = '=G19'
$ws.Range("G19").Formula = '=H19&P19'
$ws.Range("H19").Formula = '=H12'
$ws.Range("I19").Value = '15 days'
$ws.Range("N19").Formula = '=N12'
$ws.Range("P19").Value = '_15d'

# ---- Block 3: rows 20-26 (2 weeks) ----
# Row 20: Delayed transition
$ws.Range("A20").Formula = '=$C$1&TEXT(N20,"0000")'
$ws.Range("B20").Formula = '=G20'
$ws.Range("G20").Formula = '=H20&P20'
$ws.Range("H20").Formula = '=H13'
$ws.Range("I20").Value = '2 weeks'
$ws.Range("N20").Formula = '=N13'
$ws.Range("P20").Value = '_2w'

# Row 21: Net Zero 2050
$ws.Range("A21").Formula = '=$C$1&TEXT(N21,"0000")'
$ws.Range("B21").Formula = '=G21'
$ws.Range("G21").Formula = '=H21&P21'
$ws.Range("H21").Formula = '=H14'
$ws.Range("I21").Value = '2 weeks'
$ws.Range("N21").Formula = '=N14'
$ws.Range("P21").Value = '_2w'

# Row 22: NDCs
$ws.Range("A22").Formula = '=$C$1&TEXT(N22,"0000")'
$ws.Range("B22").Formula = '=G22'
$ws.Range("G22").Formula = '=H22&P22'
$ws.Range("H22").Formula = '=H15'
$ws.Range("I22").Value = '2 weeks'
$ws.Range("N22").Formula = '=N15'
$ws.Range("P22").Value = '_2w'

# Row 23: Below 2deg
$ws.Range("A23").Formula = '=$C$1&TEXT(N23,"0000")'
$ws.Range("B23").Formula = '=G23'
$ws.Range("G23").Formula = '=H23&P23'
$ws.Range("H23").Formula = '=H16'
$ws.Range("I23").Value = '2 weeks'
$ws.Range("N23").Formula = '=N16'
$ws.Range("P23").Value = '_2w'

# Row 24: Current Policies
$ws.Range("A24").Formula = '=$C$1&TEXT(N24,"0000")'
$ws.Range("B24").Formula = '=G24'
$ws.Range("G24").Formula = '=H24&P24'
$ws.Range("H24").Formula = '=H17'
$ws.Range("I24").Value = '2 weeks'
$ws.Range("N24").Formula = '=N17'
$ws.Range("P24").Value = '_2w'

# Row 25: Low demand
$ws.Range("A25").Formula = '=$C$1&TEXT(N25,"0000")'
$ws.Range("B25").Formula = '=G25'
$ws.Range("G25").Formula = '=H25&P25'
$ws.Range("H25").Formula = '=H18'
$ws.Range("I25").Value = '2 weeks'
$ws.Range("N25").Formula = '=N18'
$ws.Range("P25").Value = '_2w'

# Row 26: Fragmented World
$ws.Range("A26").Formula = '=$C$1&TEXT(N26,"0000")'
$ws.Range("B26").Formula = '=G26'
$ws.Range("G26").Formula = '=H26&P26'
$ws.Range("H26").Formula = '=H19'
$ws.Range("I26").Value = '2 weeks'
$ws.Range("N26").Formula = '=N19'
$ws.Range("P26").Value = '_2w'
